$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ------------------------------------------------------------------
# 1) Copy number formats down from row 55 onto the new rows 56:62
#    (F:G use the time format, N:U reuse their respective formats)
# ------------------------------------------------------------------
$ws.Range("F55:G55").Copy() | Out-Null
$ws.Range("F56:G62").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("N55:U55").Copy() | Out-Null
$ws.Range("N56:U62").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Raw data for the seven new measurement rows
#    (plus a Note that was also added to the last pre-existing row)
# ------------------------------------------------------------------

# Row 55 (pre-existing row gains a Note)
$ws.Cells.Item(55,11).Value = "66.25 Valve not connected"

# Row 56
$ws.Cells.Item(56,1).Value = 53
$ws.Cells.Item(56,4).Value = 0
$ws.Cells.Item(56,5).Value = 0
$ws.Cells.Item(56,6).Value = 45923.724999999999
$ws.Cells.Item(56,7).Value = 45925.472916666666
$ws.Cells.Item(56,8).Value = 8.5
$ws.Cells.Item(56,9).Value = 4.3
$ws.Cells.Item(56,10).Value = 2.5
$ws.Cells.Item(56,11).Value = "66.25 Valve not connected"

# Row 57
$ws.Cells.Item(57,1).Value = 54
$ws.Cells.Item(57,4).Value = 0
$ws.Cells.Item(57,5).Value = 0
$ws.Cells.Item(57,6).Value = 45923.724999999999
$ws.Cells.Item(57,7).Value = 45925.693749999999
$ws.Cells.Item(57,8).Value = 8.5
$ws.Cells.Item(57,9).Value = 3.75
$ws.Cells.Item(57,10).Value = 2.5
$ws.Cells.Item(57,11).Value = "66.25 Valve not connected"

# Row 58
$ws.Cells.Item(58,1).Value = 55
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(58,6).Value = 45925.716666666667
$ws.Cells.Item(58,7).Value = 45926.490277777775
$ws.Cells.Item(58,8).Value = 7.6
$ws.Cells.Item(58,9).Value = 5.75
$ws.Cells.Item(58,10).Value = 2.5
$ws.Cells.Item(58,11).Value = "68.46 Valve not connected"

# Row 59
$ws.Cells.Item(59,1).Value = 56
$ws.Cells.Item(59,2).Value = 286
$ws.Cells.Item(59,3).Value = 159
$ws.Cells.Item(59,4).Value = 0
$ws.Cells.Item(59,5).Value = 0
$ws.Cells.Item(59,6).Value = 45925.719444444447
$ws.Cells.Item(59,7).Value = 45926.490277777775
$ws.Cells.Item(59,8).Value = 12
$ws.Cells.Item(59,9).Value = 10.3
$ws.Cells.Item(59,10).Value = 2.5
$ws.Cells.Item(59,11).Value = "68.46 Valve not connected"

# Row 60
$ws.Cells.Item(60,1).Value = 57
$ws.Cells.Item(60,2).Value = 286
$ws.Cells.Item(60,3).Value = 159
$ws.Cells.Item(60,4).Value = 0
$ws.Cells.Item(60,5).Value = 0
$ws.Cells.Item(60,6).Value = 45926.496527777781
$ws.Cells.Item(60,7).Value = 45926.574999999997
$ws.Cells.Item(60,8).Value = 10.3
$ws.Cells.Item(60,9).Value = 10.1
$ws.Cells.Item(60,10).Value = 2.5
$ws.Cells.Item(60,11).Value = "68.46 Valve not connected"

# Row 61
$ws.Cells.Item(61,1).Value = 58
$ws.Cells.Item(61,4).Value = 0
$ws.Cells.Item(61,5).Value = 0
$ws.Cells.Item(61,6).Value = 45926.494444444441
$ws.Cells.Item(61,7).Value = 45926.575694444444
$ws.Cells.Item(61,8).Value = 5.75
$ws.Cells.Item(61,9).Value = 5.55
$ws.Cells.Item(61,10).Value = 2.5
$ws.Cells.Item(61,11).Value = "68.46 Valve not connected"

# Row 62
$ws.Cells.Item(62,1).Value = 59
$ws.Cells.Item(62,4).Value = 0
$ws.Cells.Item(62,5).Value = 0
$ws.Cells.Item(62,6).Value = 45926.575694444444
$ws.Cells.Item(62,7).Value = 45926.617361111108
$ws.Cells.Item(62,8).Value = 5.6
$ws.Cells.Item(62,9).Value = 4.5
$ws.Cells.Item(62,10).Value = 25
$ws.Cells.Item(62,11).Value = "68.46 Valve not connected"

# ------------------------------------------------------------------
# 3) Formulas for columns N:U, filled in the same batches the author
#    used (matching how the shared-formula groups line up in rows
#    56 / 57 / 58-59 / 60-61 / 62)
# ------------------------------------------------------------------
$pairs = @(
    @(56,56),
    @(57,57),
    @(58,59),
    @(60,61),
    @(62,62)
)

foreach ($p in $pairs) {
    $r1 = $p[0]; $r2 = $p[1]

    $ws.Range("N${r1}:N${r2}").Formula = "=G$r1-F$r1"
    $ws.Range("O${r1}:O${r2}").Formula = "=N$r1"
    $ws.Range("P${r1}:P${r2}").Formula = "=H$r1-I$r1"
    $ws.Range("Q${r1}:Q${r2}").Formula = "=ABS((E$r1-D$r1)/0.9982)"
    $ws.Range("R${r1}:R${r2}").Formula = "=J$r1*O$r1"
    $ws.Range("S${r1}:S${r2}").Formula = "=(1-ABS(Q$r1-R$r1)/R$r1)*100"
    $ws.Range("U${r1}:U${r2}").Formula = "=(1-ABS(T$r1-J$r1)/J$r1)*100%"
}

# T column uses per-row hand-entered divisors (not shared between rows)
$ws.Range("T56").Formula = "=P56*1440/2517"
$ws.Range("T57").Formula = "=P57*1440/2835"
$ws.Range("T58").Formula = "=P58*1440/1114"
$ws.Range("T59").Formula = "=P59*1440/1110"
$ws.Range("T60").Formula = "=P60*1440/113"
$ws.Range("T61").Formula = "=P61*1440/117"
$ws.Range("T62").Formula = "=P62*1440/60"

# ------------------------------------------------------------------
# 4) Grow Table1 to cover the new rows (also extends the AutoFilter)
# ------------------------------------------------------------------
$ws.ListObjects.Item("Table1").Resize($ws.Range("A3:K62"))

# ------------------------------------------------------------------
# 5) Conditional formatting ranges need to grow to row 62 as well
# ------------------------------------------------------------------
$ws.Range("U4:U55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("U4:U62"))
$ws.Range("O4:O55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("O4:O62"))

# ------------------------------------------------------------------
# 6) Restore the view: scrolled so row 16 is at the top, with Q65
#    as the active selection
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q65").Select()
